$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 6848
$ws.Range("B2").Value = "Sr. Otávio Castro"
$ws.Range("C2").Value = "Vendas"
$ws.Range("D2").Value = "Doenca"
$ws.Range("E2").Value = 8
$ws.Range("F2").Value = 45091
$ws.Range("G2").Value = 5206.94

# Row 3
$ws.Range("A3").Value = 39037
$ws.Range("B3").Value = "Paulo Teixeira"
$ws.Range("C3").Value = "Juridico"
$ws.Range("D3").Value = "Doenca"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 45087
$ws.Range("G3").Value = 3563.99

# Row 4
$ws.Range("A4").Value = 24888
$ws.Range("B4").Value = "Isaac Pastor"
$ws.Range("C4").Value = "Marketing"
$ws.Range("E4").Value = 4
$ws.Range("F4").Value = 45093
$ws.Range("G4").Value = 3990.64

# Row 5
$ws.Range("A5").Value = 55087
$ws.Range("B5").Value = "Maria Sales"
$ws.Range("C5").Value = "Recursos Humanos"
$ws.Range("D5").Value = "Problemas pessoais"
$ws.Range("E5").Value = 7
$ws.Range("F5").Value = 45088
$ws.Range("G5").Value = 7950.31

# Row 6
$ws.Range("A6").Value = 73715
$ws.Range("B6").Value = "Dr. Arthur Ferreira"
$ws.Range("C6").Value = "Recursos Humanos"
$ws.Range("D6").Value = "Doenca"
$ws.Range("E6").Value = 5
$ws.Range("F6").Value = 45104
$ws.Range("G6").Value = 4569.07

# Row 7
$ws.Range("A7").Value = 78660
$ws.Range("B7").Value = "Dr. Gustavo Henrique Cavalcante"
$ws.Range("C7").Value = "Financeiro"
$ws.Range("D7").Value = "Consulta medica"
$ws.Range("E7").Value = 6
$ws.Range("F7").Value = 45093
$ws.Range("G7").Value = 5460.19

# Row 8
$ws.Range("A8").Value = 46082
$ws.Range("B8").Value = "Olívia Carvalho"
$ws.Range("C8").Value = "Operacoes"
$ws.Range("D8").Value = "Outros"
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 45086
$ws.Range("G8").Value = 6326.42

# Row 9
$ws.Range("A9").Value = 51293
$ws.Range("B9").Value = "Yan Teixeira"
$ws.Range("C9").Value = "TI"
$ws.Range("D9").Value = "Problemas pessoais"
$ws.Range("E9").Value = 7
$ws.Range("F9").Value = 45095
$ws.Range("G9").Value = 3975.25

# Row 10
$ws.Range("A10").Value = 74886
$ws.Range("B10").Value = "Hadassa Sales"
$ws.Range("C10").Value = "Marketing"
$ws.Range("D10").Value = "Doenca"
$ws.Range("E10").Value = 1
$ws.Range("F10").Value = 45084
$ws.Range("G10").Value = 4216.69

# Row 11
$ws.Range("A11").Value = 15649
$ws.Range("B11").Value = "Esther Costela"
$ws.Range("C11").Value = "Financeiro"
$ws.Range("D11").Value = "Outros"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 45082
$ws.Range("G11").Value = 3923.09
